$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.250.24"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "'3.913.41"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'487.60"
$ws.Range("E5").Value = "  +3.59%  "
$ws.Range("D6").Value = "'146.80"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").Value = "'0.0000345"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "'43.12"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "'10.87"
$ws.Range("E13").Value = "  +4.94%  "
$ws.Range("D14").Value = "'4.540.16"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "'3.923.11"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "'14.21"
$ws.Range("E16").Value = "  -5.73%  "
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "'19.90"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'1.14"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "'68.376.28"
$ws.Range("D21").Value = "'434.82"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'3.54"
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("D23").Value = "'14.85"
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("D24").Value = "'87.81"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "'11.36"
$ws.Range("E25").Value = "  +14.68%  "
$ws.Range("D26").Value = "'11.26"
$ws.Range("E26").Value = "  +11.17%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'38.14"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'724.67"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'13.79"
$ws.Range("E31").Value = "  +2.78%  "
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("D34").Value = "'6.28"
$ws.Range("E34").Value = "  +17.73%  "
$ws.Range("D35").Value = "'41.62"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "'0.0₃0874"
$ws.Range("E36").Value = "  +5.59%  "
$ws.Range("D37").Value = "'60.36"
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("D38").Value = "'0.403"
$ws.Range("E38").Value = "  +20.11%  "
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +16.23%  "
$ws.Range("D42").Value = "'0.0482"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("D44").Value = "'2.92"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").Value = "'3.42"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").Value = "'2.14"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").Value = "'0.0₆0349"
$ws.Range("E50").Value = "  +35.79%  "
$ws.Range("D51").Value = "'144.61"
$ws.Range("E51").Value = "  -2.32%  "
